$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-15 Monday" "2024-07-16 Tuesday"

Replace-Text "676×5=3380" "239×9=2151"
Replace-Text "274×6=1644" "359×9=3231"
Replace-Text "674×6=4044" "889×2=1778"
Replace-Text "851×2=1702" "899×2=1798"
Replace-Text "908×6=5448" "889×5=4445"
Replace-Text "509×2=1018" "390×9=3510"
Replace-Text "667×6=4002" "390×7=2730"
Replace-Text "171×6=1026" "939×7=6573"
Replace-Text "338×3=1014" "170×8=1360"
Replace-Text "312×3=936" "120×6=720"
Replace-Text "480×3=1440" "233×6=1398"
Replace-Text "623×3=1869" "318×7=2226"
Replace-Text "268×9=2412" "208×3=624"
Replace-Text "586×7=4102" "792×7=5544"
Replace-Text "206×5=1030" "257×9=2313"
Replace-Text "295×6=1770" "871×6=5226"
Replace-Text "461×7=3227" "496×8=3968"
Replace-Text "272×2=544" "730×5=3650"
Replace-Text "895×2=1790" "367×9=3303"
Replace-Text "634×2=1268" "740×4=2960"
Replace-Text "846×5=4230" "436×2=872"
Replace-Text "710×3=2130" "328×8=2624"
Replace-Text "770×3=2310" "602×5=3010"
Replace-Text "912×3=2736" "494×2=988"
Replace-Text "725×7=5075" "418×4=1672"
